# Update "想去人数" (F column) counts on the 展览, 演出, and 全部类型 sheets
# to reflect the refreshed data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 3358
    4  = 130
    5  = 6947
    6  = 2350
    8  = 98
    13 = 170
    14 = 560
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# 演出 sheet
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
    2 = 21
}
foreach ($row in $updates2.Keys) {
    $ws2.Range("F$row").Value = $updates2[$row]
}

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 3358
    3  = 21
    5  = 130
    6  = 6947
    7  = 2350
    9  = 98
    14 = 170
    15 = 560
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
